$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to be treated as text so that
# numeric-looking values (e.g. "9.19", "1.00") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

# --- Simple value updates (rows where only D and/or E change) ---
$ws.Range("D2").Value = "44.113.76"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "2.362.97"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "0.698"
$ws.Range("E5").Value = "  +6.26%  "
$ws.Range("D6").Value = "241.91"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("D7").Value = "76.54"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +18.87%  "
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("D11").Value = "57.38"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "33.44"
$ws.Range("E12").Value = "  +21.78%  "
$ws.Range("D13").Value = "7.48"
$ws.Range("E13").Value = "  +12.03%  "
$ws.Range("D14").Value = "0.108"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "2.714.35"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "16.72"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "0.929"
$ws.Range("E17").Value = "  +5.41%  "
$ws.Range("D18").Value = "2.360.83"
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("D19").Value = "43.982.34"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  +5.92%  "
$ws.Range("D22").Value = "77.92"
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("D23").Value = "260.94"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").Value = "1.81"
$ws.Range("E27").Value = "  +17.56%  "
$ws.Range("D28").Value = "10.97"
$ws.Range("E28").Value = "  +6.86%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "23.26"
$ws.Range("E30").Value = "  +3.82%  "
$ws.Range("D31").Value = "174.95"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  +6.42%  "
$ws.Range("D35").Value = "0.0764"
$ws.Range("E35").Value = "  +9.38%  "
$ws.Range("D36").Value = "5.46"
$ws.Range("E36").Value = "  +7.25%  "
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").Value = "2.44"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "6.43"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +7.32%  "
$ws.Range("D41").Value = "0.221"
$ws.Range("E41").Value = "  +23.43%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "4.83"
$ws.Range("E46").Value = "  +8.76%  "
$ws.Range("E47").Value = "  +9.96%  "
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("D49").Value = "102.30"
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").Value = "56.47"
$ws.Range("E51").Value = "  +10.36%  "

# --- Rows 42-44 reordering (Cronos / InjectiveProtocol / FraxShare) ---
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "9.19"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.108"
$ws.Range("E43").Value = "  +12.38%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "19.29"
$ws.Range("E44").Value = "  -0.75%  "

# Reset style back to Normal so we do not leave a stray text-format style
# applied to the cells (keeps styling identical to the original workbook).
$priceRange.Style = "Normal"
$volRange.Style = "Normal"
